$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-08-15 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-16 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("72÷4=18, 0", $true, $false, $false, $false, $false, $true, 1, $false, "86÷7=12, 2", 2) | Out-Null
$d.Content.Find.Execute("65÷5=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "25÷3=8, 1", 2) | Out-Null
$d.Content.Find.Execute("35÷9=3, 8", $true, $false, $false, $false, $false, $true, 1, $false, "85÷8=10, 5", 2) | Out-Null
$d.Content.Find.Execute("83÷2=41, 1", $true, $false, $false, $false, $false, $true, 1, $false, "51÷3=17, 0", 2) | Out-Null
$d.Content.Find.Execute("33÷4=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "99÷2=49, 1", 2) | Out-Null
$d.Content.Find.Execute("30÷7=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "19÷7=2, 5", 2) | Out-Null
$d.Content.Find.Execute("37÷4=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "70÷2=35, 0", 2) | Out-Null
$d.Content.Find.Execute("59÷2=29, 1", $true, $false, $false, $false, $false, $true, 1, $false, "80÷4=20, 0", 2) | Out-Null
$d.Content.Find.Execute("42÷5=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "49÷3=16, 1", 2) | Out-Null
$d.Content.Find.Execute("84÷8=10, 4", $true, $false, $false, $false, $false, $true, 1, $false, "51÷4=12, 3", 2) | Out-Null
$d.Content.Find.Execute("72÷5=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "69÷7=9, 6", 2) | Out-Null
$d.Content.Find.Execute("30÷5=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "61÷4=15, 1", 2) | Out-Null
$d.Content.Find.Execute("38÷7=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "57÷8=7, 1", 2) | Out-Null
$d.Content.Find.Execute("58÷9=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "45÷8=5, 5", 2) | Out-Null
$d.Content.Find.Execute("49÷4=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "54÷4=13, 2", 2) | Out-Null
$d.Content.Find.Execute("20÷2=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "47÷4=11, 3", 2) | Out-Null
$d.Content.Find.Execute("60÷9=6, 6", $true, $false, $false, $false, $false, $true, 1, $false, "16÷6=2, 4", 2) | Out-Null
$d.Content.Find.Execute("99÷6=16, 3", $true, $false, $false, $false, $false, $true, 1, $false, "95÷4=23, 3", 2) | Out-Null
$d.Content.Find.Execute("73÷4=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "55÷4=13, 3", 2) | Out-Null
$d.Content.Find.Execute("45÷3=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "53÷3=17, 2", 2) | Out-Null
$d.Content.Find.Execute("50÷2=25, 0", $true, $false, $false, $false, $false, $true, 1, $false, "21÷8=2, 5", 2) | Out-Null
$d.Content.Find.Execute("54÷9=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "22÷5=4, 2", 2) | Out-Null
$d.Content.Find.Execute("76÷6=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "21÷6=3, 3", 2) | Out-Null
$d.Content.Find.Execute("76÷2=38, 0", $true, $false, $false, $false, $false, $true, 1, $false, "47÷8=5, 7", 2) | Out-Null
$d.Content.Find.Execute("92÷5=18, 2", $true, $false, $false, $false, $false, $true, 1, $false, "79÷6=13, 1", 2) | Out-Null
